$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp value in A13 (floating point precision fix)
$ws.Range("A13").Value = 45876.4585312037

# Append new row 14 with the latest sensor reading
$ws.Range("A14").Value = 45876.50018915661
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = 28
$ws.Range("D14").Value = 19.06
$ws.Range("E14").Value = 78.45999999999999
$ws.Range("F14").Value = 631.05
$ws.Range("G14").Value = 12.16
$ws.Range("H14").Value = "ESE"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "12:00:16"

# Match the number format of the A column datetime cells (style index 2 in the original sheet)
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
